# "GOOD DATA FOR 8 GANTT CHARTS"
# Corrects the treatment / year / length data for rows 22-24 on the
# "Export Worksheet" sheet, sets explicit (best-fit) column widths for
# columns A:I, and moves the active cell selection to E25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: treatment corrected to "Digouts" -------------------------
$ws.Range("B22").Value = "Digouts"

# --- Row 23: EA / treatment / year / length / budget group corrected --
$ws.Range("A23").Value = "04-2E330"
$ws.Range("B23").Value = "Digouts"
$ws.Range("E23").Value = 2011
$ws.Range("F23").Value = 12
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = "HM Corrective"

# --- Row 24: EA / treatment / length / budget group corrected ---------
$ws.Range("A24").Value = "04-2E510"
$ws.Range("B24").Value = "HMA Thin Overlay"
$ws.Range("F24").Value = 0
$ws.Range("H24").Value = 13
$ws.Range("I24").Value = "HM Preventive"

# --- Explicit (best-fit) column widths for A:I -------------------------
# (values chosen so the saved column width, after this runtime's
# internal rounding, lands as close as possible to the target widths
# 10.140625, 22, 8.28515625, 6.85546875, 5.5703125, 8.28515625,
# 8.5703125, 7.85546875, 15.42578125)
$ws.Columns.Item(1).ColumnWidth = 9.333333333333334
$ws.Columns.Item(2).ColumnWidth = 21.166666666666668
$ws.Columns.Item(3).ColumnWidth = 7.5
$ws.Columns.Item(4).ColumnWidth = 6
$ws.Columns.Item(5).ColumnWidth = 4.666666666666667
$ws.Columns.Item(6).ColumnWidth = 7.5
$ws.Columns.Item(7).ColumnWidth = 7.666666666666667
$ws.Columns.Item(8).ColumnWidth = 7
$ws.Columns.Item(9).ColumnWidth = 14.666666666666666

# --- Move the active selection from E24 to E25 -------------------------
$ws.Range("E25").Select()
